# Grade update: fill in remaining scores/letters for Loren Grey, Esiete Yismaw
# Mebratie, Dary Ductoc and Penelope Turgen.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Will Gearhart: grade corrected from F to D+
$ws.Range("M2").Value = "D+"

# Row 6 - Loren Grey: letter grade comes in ("??" - pending review)
$ws.Range("M6").Value = "??"

# Row 7 - Esiete Yismaw Mebratie: letter grade comes in
$ws.Range("M7").Value = "B-"

# Row 8 - Dary Ductoc: Midterm 1 score entered, letter grade comes in
$ws.Range("F8").Formula = "=43/50"
$ws.Range("M8").Value = "A"

# Row 10 - Penelope Turgen: Science Paper, Journal Reflection and
# In-class participation scores entered, letter grade comes in
$ws.Range("G10").Formula = "=50/60"
$ws.Range("H10").Formula = "=40/100"
$ws.Range("I10").Formula = "=1"
$ws.Range("M10").Value = "??"

$ws.Range("I11").Select()
